# Dufour1992.xlsx - "Healthy status added to experimental data"
#
# Content changes:
#   - shared string "state"  -> "status"   (column C header, rows 3 & 4)
#   - shared string "normal" -> "healthy"  (column C data, rows 5-10)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "state" -> "status"
$ws.Range("C3").Value = "status"
$ws.Range("C4").Value = "status"

# Data column: "normal" -> "healthy"
$ws.Range("C5").Value = "healthy"
$ws.Range("C6").Value = "healthy"
$ws.Range("C7").Value = "healthy"
$ws.Range("C8").Value = "healthy"
$ws.Range("C9").Value = "healthy"
$ws.Range("C10").Value = "healthy"

# Minor layout adjustments accompanying the edit session
$ws.Rows.Item(6).RowHeight = 12.8
$ws.Rows.Item(7).RowHeight = 12.8
$ws.Rows.Item(8).RowHeight = 12.8
$ws.Rows.Item(9).RowHeight = 12.8
$ws.Rows.Item(10).RowHeight = 12.8

$ws.Columns.Item(8).ColumnWidth = 11.666666666666666

$ws.Range("A4:I19").Select() | Out-Null
